$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update existing row 15 (Dia 14, July) total_venda value
$ws.Range("B15").Value = 21313.41

# 2. Insert a new row at position 16 (shifts everything below down by one, all the
#    way to the end of the table), then fill it with the new daily record
#    (Dia 15, July, value 32100.04)
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 32100.04
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 2025
$ws.Range("E16").Value = "07/2025"
